$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.842.01"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "'3.763.40"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'426.93"
$ws.Range("E5").Value = "  +5.37%  "
$ws.Range("D6").Value = "'139.09"
$ws.Range("E6").Value = "  +5.49%  "
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  +2.69%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.726"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'0.149"
$ws.Range("E10").Value = "  -10.66%  "
$ws.Range("D11").Value = "'0.0000305"
$ws.Range("E11").Value = "  -16.46%  "
$ws.Range("D12").Value = "'42.36"
$ws.Range("E12").Value = "  +4.04%  "
$ws.Range("E13").Value = "  +5.60%  "
$ws.Range("D14").Value = "'4.371.71"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").Value = "'14.89"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "'3.783.94"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "'19.78"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("E19").Value = "  +5.36%  "
$ws.Range("D20").Value = "'65.981.28"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").Value = "'402.46"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").Value = "'14.78"
$ws.Range("E22").Value = "  +3.13%  "
$ws.Range("D23").Value = "'3.28"
$ws.Range("E23").Value = "  +7.63%  "
$ws.Range("D24").Value = "'84.27"
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("D25").Value = "'10.16"
$ws.Range("E25").Value = "  +38.56%  "
$ws.Range("D26").Value = "'36.36"
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("E27").Value = "  +4.95%  "
$ws.Range("E28").Value = "  -3.24%  "
$ws.Range("D29").Value = "'9.72"
$ws.Range("E29").Value = "  +4.14%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.134"
$ws.Range("E30").Value = "  +12.05%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").Value = "'13.59"
$ws.Range("E31").Value = "  +10.18%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'699.70"
$ws.Range("E32").Value = "  -5.12%  "
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("D34").Value = "'40.82"
$ws.Range("E34").Value = "  +5.06%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = "'5.75"
$ws.Range("E36").Value = "  +35.35%  "
$ws.Range("E37").Value = "  -3.65%  "
$ws.Range("D38").Value = "'56.08"
$ws.Range("E38").Value = "  +2.27%  "
$ws.Range("D39").Value = "'0.0467"
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("D40").Value = "'2.79"
$ws.Range("E40").Value = "  +40.11%  "
$ws.Range("E41").Value = "  +3.68%  "
$ws.Range("E42").Value = "  +4.33%  "
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("D44").Value = "'0.0₃0659"
$ws.Range("E44").Value = "  -12.35%  "
$ws.Range("D45").Value = "'0.324"
$ws.Range("E45").Value = "  +10.72%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.19"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").Value = "'3.34"
$ws.Range("E47").Value = "  +2.89%  "
$ws.Range("D48").Value = "'2.66"
$ws.Range("E48").Value = "  +2.13%  "
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").Value = "'138.69"
$ws.Range("E50").Value = "  -4.39%  "
$ws.Range("E51").Value = "  -0.41%  "
